# Auto-generated edit script: adds "ECs" as a new sending cluster row-block
# and re-derives the M1/M2 Il1a->Il1r1 edge metrics (Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Il1a -> Il1r1 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1a"
$ws.Range("C2").Value = "Il1r1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02286966666666667
$ws.Range("H2").Value = 0.068609
$ws.Range("I2").Value = 0.001711767187487096
$ws.Range("J2").Value = 0.001711767187487096
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.37432
$ws.Range("N2").Value = 40.12296
$ws.Range("O2").Value = 0.1019828318366699
$ws.Range("P2").Value = 0.1077302617359031
$ws.Range("Q2").Value = 0.3058662402933333
$ws.Range("R2").Value = 2.75279616264
$ws.Range("S2").Value = 0.0001745708652250259
$ws.Range("T2").Value = 0.0001844091271389156

# Row 3: ECs -> Il1a -> Il1r1 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il1a"
$ws.Range("C3").Value = "Il1r1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02286966666666667
$ws.Range("H3").Value = 0.068609
$ws.Range("I3").Value = 0.001711767187487096
$ws.Range("J3").Value = 0.001711767187487096
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 96.23965333333335
$ws.Range("N3").Value = 288.71896
$ws.Range("O3").Value = 0.7338535627914345
$ws.Range("P3").Value = 0.7752112289052887
$ws.Range("Q3").Value = 2.200968791848889
$ws.Range("R3").Value = 19.80871912664
$ws.Range("S3").Value = 0.001256186449206879
$ws.Range("T3").Value = 0.001326981145011622

# Row 4: ECs -> Il1a -> Il1r1 -> M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il1a"
$ws.Range("C4").Value = "Il1r1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02286966666666667
$ws.Range("H4").Value = 0.068609
$ws.Range("I4").Value = 0.001711767187487096
$ws.Range("J4").Value = 0.001711767187487096
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1994553333333333
$ws.Range("N4").Value = 0.598366
$ws.Range("O4").Value = 0.001520901228493132
$ws.Range("P4").Value = 0.001606614412143705
$ws.Range("Q4").Value = 0.004561476988222222
$ws.Range("R4").Value = 0.041053292894
$ws.Range("S4").Value = 0.000002603428818343358
$ws.Range("T4").Value = 0.000002750149833651465

# Row 5: ECs -> Il1a -> Il1r1 -> M2
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il1a"
$ws.Range("C5").Value = "Il1r1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02286966666666667
$ws.Range("H5").Value = 0.068609
$ws.Range("I5").Value = 0.001711767187487096
$ws.Range("J5").Value = 0.001711767187487096
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3399400000000001
$ws.Range("N5").Value = 1.01982
$ws.Range("O5").Value = 0.002592135065899243
$ws.Range("P5").Value = 0.002738219601034139
$ws.Range("Q5").Value = 0.007774314486666668
$ws.Range("R5").Value = 0.06996883038000001
$ws.Range("S5").Value = 0.000004437131751341026
$ws.Range("T5").Value = 0.000004687194465184248

# Row 6: ECs -> Il1a -> Il1r1 -> sCs
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Il1a"
$ws.Range("C6").Value = "Il1r1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02286966666666667
$ws.Range("H6").Value = 0.068609
$ws.Range("I6").Value = 0.001711767187487096
$ws.Range("J6").Value = 0.001711767187487096
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 20.989489
$ws.Range("N6").Value = 41.978978
$ws.Range("O6").Value = 0.1600505690775031
$ws.Range("P6").Value = 0.1127136753456305
$ws.Range("Q6").Value = 0.4800226169336667
$ws.Range("R6").Value = 2.880135701602
$ws.Range("S6").Value = 0.0002739693124855068
$ws.Range("T6").Value = 0.0001929395710377236

# Row 7: M1 -> Il1a -> Il1r1 -> ECs
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Il1a"
$ws.Range("C7").Value = "Il1r1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.046308333333333
$ws.Range("H7").Value = 15.138925
$ws.Range("I7").Value = 0.3777101410722805
$ws.Range("J7").Value = 0.3777101410722805
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.37432
$ws.Range("N7").Value = 40.12296
$ws.Range("O7").Value = 0.1019828318366699
$ws.Range("P7").Value = 0.1077302617359031
$ws.Range("Q7").Value = 67.49094246866666
$ws.Range("R7").Value = 607.418482218
$ws.Range("S7").Value = 0.03851994979997923
$ws.Range("T7").Value = 0.04069081235802165

# Row 8: M1 -> Il1a -> Il1r1 -> FAPs
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Il1a"
$ws.Range("C8").Value = "Il1r1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.046308333333333
$ws.Range("H8").Value = 15.138925
$ws.Range("I8").Value = 0.3777101410722805
$ws.Range("J8").Value = 0.3777101410722805
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 96.23965333333335
$ws.Range("N8").Value = 288.71896
$ws.Range("O8").Value = 0.7338535627914345
$ws.Range("P8").Value = 0.7752112289052887
$ws.Range("Q8").Value = 485.6549646131112
$ws.Range("R8").Value = 4370.894681518001
$ws.Range("S8").Value = 0.2771839327283483
$ws.Range("T8").Value = 0.2928051426306325

# Row 9: M1 -> Il1a -> Il1r1 -> M1
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Il1a"
$ws.Range("C9").Value = "Il1r1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.046308333333333
$ws.Range("H9").Value = 15.138925
$ws.Range("I9").Value = 0.3777101410722805
$ws.Range("J9").Value = 0.3777101410722805
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1994553333333333
$ws.Range("N9").Value = 0.598366
$ws.Range("O9").Value = 0.001520901228493132
$ws.Range("P9").Value = 0.001606614412143705
$ws.Range("Q9").Value = 1.006513110727778
$ws.Range("R9").Value = 9.05861799655
$ws.Range("S9").Value = 0.0005744598175711456
$ws.Range("T9").Value = 0.0006068345562595578

# Row 10: M1 -> Il1a -> Il1r1 -> M2
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Il1a"
$ws.Range("C10").Value = "Il1r1"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.046308333333333
$ws.Range("H10").Value = 15.138925
$ws.Range("I10").Value = 0.3777101410722805
$ws.Range("J10").Value = 0.3777101410722805
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3399400000000001
$ws.Range("N10").Value = 1.01982
$ws.Range("O10").Value = 0.002592135065899243
$ws.Range("P10").Value = 0.002738219601034139
$ws.Range("Q10").Value = 1.715442054833334
$ws.Range("R10").Value = 15.4389784935
$ws.Range("S10").Value = 0.000979075701419208
$ws.Range("T10").Value = 0.001034253311793488

# Row 11: M1 -> Il1a -> Il1r1 -> sCs
$ws.Range("A11").Value = "M1"
$ws.Range("B11").Value = "Il1a"
$ws.Range("C11").Value = "Il1r1"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5.046308333333333
$ws.Range("H11").Value = 15.138925
$ws.Range("I11").Value = 0.3777101410722805
$ws.Range("J11").Value = 0.3777101410722805
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 20.989489
$ws.Range("N11").Value = 41.978978
$ws.Range("O11").Value = 0.1600505690775031
$ws.Range("P11").Value = 0.1127136753456305
$ws.Range("Q11").Value = 105.9194332531083
$ws.Range("R11").Value = 635.51659951865
$ws.Range("S11").Value = 0.06045272302496248
$ws.Range("T11").Value = 0.04257309821557331

# Row 12: M2 -> Il1a -> Il1r1 -> ECs
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Il1a"
$ws.Range("C12").Value = "Il1r1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 8.291089
$ws.Range("H12").Value = 24.873267
$ws.Range("I12").Value = 0.6205780917402324
$ws.Range("J12").Value = 0.6205780917402324
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 13.37432
$ws.Range("N12").Value = 40.12296
$ws.Range("O12").Value = 0.1019828318366699
$ws.Range("P12").Value = 0.1077302617359031
$ws.Range("Q12").Value = 110.88767743448
$ws.Range("R12").Value = 997.9890969103199
$ws.Range("S12").Value = 0.06328831117146559
$ws.Range("T12").Value = 0.06685504025074251

# Row 13: M2 -> Il1a -> Il1r1 -> FAPs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Il1a"
$ws.Range("C13").Value = "Il1r1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 8.291089
$ws.Range("H13").Value = 24.873267
$ws.Range("I13").Value = 0.6205780917402324
$ws.Range("J13").Value = 0.6205780917402324
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 96.23965333333335
$ws.Range("N13").Value = 288.71896
$ws.Range("O13").Value = 0.7338535627914345
$ws.Range("P13").Value = 0.7752112289052887
$ws.Range("Q13").Value = 797.9315311158134
$ws.Range("R13").Value = 7181.383780042321
$ws.Range("S13").Value = 0.4554134436138793
$ws.Range("T13").Value = 0.4810791051296445

# Row 14: M2 -> Il1a -> Il1r1 -> M1
$ws.Range("A14").Value = "M2"
$ws.Range("B14").Value = "Il1a"
$ws.Range("C14").Value = "Il1r1"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 8.291089
$ws.Range("H14").Value = 24.873267
$ws.Range("I14").Value = 0.6205780917402324
$ws.Range("J14").Value = 0.6205780917402324
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.1994553333333333
$ws.Range("N14").Value = 0.598366
$ws.Range("O14").Value = 0.001520901228493132
$ws.Range("P14").Value = 0.001606614412143705
$ws.Range("Q14").Value = 1.653701920191333
$ws.Range("R14").Value = 14.883317281722
$ws.Range("S14").Value = 0.000943837982103643
$ws.Range("T14").Value = 0.0009970297060504958

# Row 15: M2 -> Il1a -> Il1r1 -> M2
$ws.Range("A15").Value = "M2"
$ws.Range("B15").Value = "Il1a"
$ws.Range("C15").Value = "Il1r1"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 8.291089
$ws.Range("H15").Value = 24.873267
$ws.Range("I15").Value = 0.6205780917402324
$ws.Range("J15").Value = 0.6205780917402324
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.3399400000000001
$ws.Range("N15").Value = 1.01982
$ws.Range("O15").Value = 0.002592135065899243
$ws.Range("P15").Value = 0.002738219601034139
$ws.Range("Q15").Value = 2.81847279466
$ws.Range("R15").Value = 25.36625515194
$ws.Range("S15").Value = 0.001608622232728694
$ws.Range("T15").Value = 0.001699279094775466

# Row 16: M2 -> Il1a -> Il1r1 -> sCs
$ws.Range("A16").Value = "M2"
$ws.Range("B16").Value = "Il1a"
$ws.Range("C16").Value = "Il1r1"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 8.291089
$ws.Range("H16").Value = 24.873267
$ws.Range("I16").Value = 0.6205780917402324
$ws.Range("J16").Value = 0.6205780917402324
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 20.989489
$ws.Range("N16").Value = 41.978978
$ws.Range("O16").Value = 0.1600505690775031
$ws.Range("P16").Value = 0.1127136753456305
$ws.Range("Q16").Value = 174.025721363521
$ws.Range("R16").Value = 1044.154328181126
$ws.Range("S16").Value = 0.09932387674005515
$ws.Range("T16").Value = 0.06994763755901943

